$d = $word.ActiveDocument

# Simple whole-text replacement (safe when the matched run either has no
# sibling empty run, or carries distinguishing run properties so the
# engine won't silently coalesce it with a neighboring run).
function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# Text-only replacement that preserves any empty sibling run (e.g. the
# leading "<w:r/>" placeholder runs used throughout this document) by
# rewriting just the matched run via InsertXML instead of doing a
# Find/Replace (which here coalesces adjacent runs that share formatting).
function Replace-RunText($old, $new) {
    $searchRng = $d.Content
    $found = $searchRng.Find.Execute($old, $true, $true, $false, $false, `
                                      $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Text not found: $old"
    }
    $rng = $d.Range($searchRng.Start, $searchRng.End)
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage" ' + `
           'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
           '<pkg:part pkg:name="/word/document.xml" ' + `
           'pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
           '<pkg:xmlData><w:document><w:body><w:p><w:r><w:t>' + $new + `
           '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($xml)
}

# Title change (appears twice: main heading + bold summary line; both
# runs carry distinguishing context so a plain replace-all is safe)
Replace-Text "Play Alpha Eagle Stack N Sync for Free - Exciting Winter Slot Game" "Play Alpha Eagle Stack N Sync for Free - Exciting Slot Game"

# "What we like" bullet list
Replace-RunText "Beautifully designed icons with a winter theme" "Stunning graphics and immersive snowy mountain background"
Replace-RunText "High volatility rate for a chance to win big" "Beautifully designed symbols with a cool winter theme"
Replace-RunText "Great RTP percentage of 96.26%" "Exciting bonus symbol system with various special features"
Replace-RunText "Variety of special features including Free Spins and Respins" "High volatility rate and excellent RTP percentage"

# "What we don't like" bullet list
Replace-RunText "No progressive jackpot feature" "Limited number of bonus symbols"
Replace-RunText "Bonus Buy button may not be suitable for all players" "Bonus features can only be purchased with the Bonus Buy button"

# Closing italic summary line
Replace-Text "Experience the exciting winter-themed slot game Alpha Eagle Stack N Sync for free. Enjoy Respins, Free Spins, and a variety of bonus features." "Read our review of Alpha Eagle Stack N Sync, an exciting slot game with stunning graphics. Play for free and experience the winter theme."
